$d = $word.ActiveDocument

# Locate the paragraph containing "In Replay Mode" so the new table can be
# inserted immediately after it (and before the following blank paragraph).
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "In Replay Mode`r") {
        $targetPara = $p
    }
}

$insertRange = $d.Range($targetPara.Range.End, $targetPara.Range.End)

$tableXml = '<w:tbl>' + `
    '<w:tblPr>' + `
        '<w:tblStyle w:val="TableGrid"/>' + `
        '<w:tblW w:w="0" w:type="auto"/>' + `
        '<w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>' + `
    '</w:tblPr>' + `
    '<w:tblGrid>' + `
        '<w:gridCol w:w="4675"/>' + `
        '<w:gridCol w:w="4675"/>' + `
    '</w:tblGrid>' + `
    '<w:tr><w:tc><w:tcPr><w:tcW w:w="4675" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>Action</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4675" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>Key</w:t></w:r></w:p></w:tc></w:tr>' + `
    '<w:tr><w:tc><w:tcPr><w:tcW w:w="4675" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>Return / Settings Menu</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4675" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>ESC</w:t></w:r></w:p></w:tc></w:tr>' + `
    '<w:tr><w:tc><w:tcPr><w:tcW w:w="4675" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>Move camera view</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4675" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>Mouse Right Click + Mouse Movement</w:t></w:r></w:p></w:tc></w:tr>' + `
    '<w:tr><w:tc><w:tcPr><w:tcW w:w="4675" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>Zoom in/out</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4675" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>Mouse scroll</w:t></w:r></w:p></w:tc></w:tr>' + `
    '<w:tr><w:tc><w:tcPr><w:tcW w:w="4675" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>Change Drilling Leader Tower details visibility</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4675" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>V</w:t></w:r></w:p></w:tc></w:tr>' + `
    '<w:tr><w:tc><w:tcPr><w:tcW w:w="4675" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>Change Terrain Layer visibility</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4675" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>T</w:t></w:r></w:p></w:tc></w:tr>' + `
    '</w:tbl>'

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + $tableXml + '</w:body>' + `
    '</w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

$insertRange.InsertXML($packageXml)
